# Weekly update: insert three new daily price rows (date 45275) at the top
# of the Chirimoya price history table (Femacal de La Calera), pushing the
# existing rows down by three positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 430:432, shifting rows 430..512 down to 433..515.
$ws.Rows("430:432").Insert()

# New data for the three inserted rows (A..T), matching the rest of the
# table's constant columns (market/product identifiers) plus the new
# week's values.
$newRows = @(
    @(3, "Femacal de La Calera", "Coquimbo", 45275, 5, "Fruta", 100107, "Otros", 100107002, "Chirimoya", "Cultivar IV Región", "Especial", 45, 22000, 22000, 22000, '$/bandeja 10 kilos', "Provincia del Elquí", 2200, 10),
    @(3, "Femacal de La Calera", "Coquimbo", 45275, 5, "Fruta", 100107, "Otros", 100107002, "Chirimoya", "Cultivar IV Región", "Primera", 56, 20000, 20000, 20000, '$/bandeja 10 kilos', "Provincia del Elquí", 2000, 10),
    @(3, "Femacal de La Calera", "Coquimbo", 45275, 5, "Fruta", 100107, "Otros", 100107002, "Chirimoya", "Cultivar IV Región", "Segunda", 45, 17000, 17000, 17000, '$/bandeja 10 kilos', "Provincia del Elquí", 1700, 10)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowNum = 430 + $i
    $rowValues = $newRows[$i]
    for ($col = 1; $col -le $rowValues.Count; $col++) {
        $ws.Cells.Item($rowNum, $col).Value2 = $rowValues[$col - 1]
    }
}
